# "validando entrada y mejorando codigo"
# Rebuild the DFA transition table on Hoja1:
#   Estados | A  | B  | Aceptador
#   s0      | s0 | s1 | si
#   s1      | s2 | s1 | no
#   s2      | s2 | s3 | si
#   s3      | B  | B  | no
# and mark F5 with an underlined (but empty) style, adjust the page setup
# and selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Start from a clean sheet (the previous "estados/a /b" table is fully replaced).
$ws.Cells.Clear()

# Fill the states column first ...
$ws.Range("A2").Value = "s0"
$ws.Range("A3").Value = "s1"
$ws.Range("A4").Value = "s2"
$ws.Range("A5").Value = "s3"

# ... then the header row ...
$ws.Range("B1").Value = "A"
$ws.Range("C1").Value = "B"
$ws.Range("A1").Value = "Estados"
$ws.Range("D1").Value = "Aceptador"

# ... then the transition data, row by row.
$ws.Range("B2").Value = "s0"
$ws.Range("C2").Value = "s1"
$ws.Range("D2").Value = "si"

$ws.Range("B3").Value = "s2"
$ws.Range("C3").Value = "s1"
$ws.Range("D3").Value = "no"

$ws.Range("B4").Value = "s2"
$ws.Range("C4").Value = "s3"
$ws.Range("D4").Value = "si"

$ws.Range("B5").Value = "B"
$ws.Range("C5").Value = "B"
$ws.Range("D5").Value = "no"

# Empty validation marker cell with an underline style applied.
$ws.Range("F5").Font.Underline = $true

# Match the saved selection / active cell.
[void]$ws.Range("B5").Select()

# Page setup tweaks (A4 portrait).
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1
